$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AdventureLevelTable")

# Update the escalating gold reward values (E6:E13) to a flat 10000
$ws.Range("E6:E13").Value = "10000"

# Update sheet view: reset scroll position and change active selection
$ws.Range("D9").Select()
